# Update heapsort results with new measured values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("heapsort")

# The two result rows swap labels: row 2 is now the AVERAGE_CASE series,
# row 3 is now the BEST_CASE series, each with freshly measured data.
$ws.Range("A2").Value = "AVERAGE_CASE"
$ws.Range("A3").Value = "BEST_CASE"

$row2 = 8160, 18020, 39220, 84860, 180460, 378480, 757480, 1508130, 3012890
$row3 = 510, 980, 1900, 3770, 7510, 14950, 29890, 59660, 119900

for ($i = 0; $i -lt 9; $i++) {
    $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
    $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
}

# Chart was moved/resized on the sheet.
$co = $ws.ChartObjects().Item(1)
$co.Left = 13.12488188976378
$co.Top = 95.62488188976378
$co.Width = 523.7949612450788
$co.Height = 288.37503937007875

# Selection left on J22.
$null = $ws.Range("J22").Select()
